$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (K1:N1) -------------------------------------------
# K1/L1/M1 are new columns; N1 takes over the old "performance_improvement"
# header that used to live in K1.
$ws.Range("K1").Value = "projected_tgs"
$ws.Range("L1").Value = "theoretical_tgs"
$ws.Range("M1").Value = "output_token_rate"
$ws.Range("N1").Value = "performance_improvement"

# Copy the header style (bold/centered/bordered) from an existing header
# cell onto the three newly-created header cells, as well as re-applying
# it to N1 (whose "performance_improvement" text moved out of K1).
$ws.Range("J1").Copy() | Out-Null
$ws.Range("K1:N1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Row 2 (600MHz) -------------------------------------------------------
$ws.Range("E2").Value = 22229.33333299483
$ws.Range("F2").Value = 142.5599999978292
$ws.Range("G2").Value = 22371.89333299266
$ws.Range("H2").Value = 22229.33333333333
$ws.Range("I2").Value = 142.56
$ws.Range("J2").Value = 22371.89333333333
$ws.Range("K2").Value = 5.095679578977786
$ws.Range("L2").Value = 5.095679578900191
$ws.Range("M2").Value = 14.02918069606099
$ws.Range("N2").Value = -62.49999999942897

# --- Row 3 (1000MHz) -------------------------------------------------------
$ws.Range("E3").Value = 13337.6
$ws.Range("F3").Value = 85.536
$ws.Range("G3").Value = 13423.136
$ws.Range("H3").Value = 13337.6
$ws.Range("I3").Value = 85.536
$ws.Range("J3").Value = 13423.136
$ws.Range("K3").Value = 8.492799298166984
$ws.Range("L3").Value = 8.492799298166984
$ws.Range("M3").Value = 23.38196782641227
$ws.Range("N3").Value = -37.50000000000001

# --- Row 4 (1600MHz) -------------------------------------------------------
$ws.Range("E4").Value = 8336
$ws.Range("F4").Value = 53.46
$ws.Range("G4").Value = 8389.459999999999
$ws.Range("H4").Value = 8336
$ws.Range("I4").Value = 53.46
$ws.Range("J4").Value = 8389.459999999999
$ws.Range("K4").Value = 13.58847887706718
$ws.Range("L4").Value = 13.58847887706718
$ws.Range("M4").Value = 37.41114852225963
$ws.Range("N4").Value = 0

# --- Row 5 (2000MHz) -------------------------------------------------------
$ws.Range("E5").Value = 6668.8
$ws.Range("F5").Value = 42.768
$ws.Range("G5").Value = 6711.568
$ws.Range("H5").Value = 6668.8
$ws.Range("I5").Value = 42.768
$ws.Range("J5").Value = 6711.568
$ws.Range("K5").Value = 16.98559859633397
$ws.Range("L5").Value = 16.98559859633397
$ws.Range("M5").Value = 46.76393565282454
$ws.Range("N5").Value = 24.99999999999998
